$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 675
$ws.Range("A675").NumberFormat = "@"
$ws.Range("A675").Value = "2024-09-02"
$ws.Range("A675").ClearFormats()
$ws.Range("C675").Value = 1840.550048828125
$ws.Range("D675").Value = 608.5800170898438
$ws.Range("E675").Value = 1111.550048828125
$ws.Range("F675").Value = 177.5399932861328
$ws.Range("G675").Value = 670.2000122070312
$ws.Range("H675").Value = 22669.95040893555
$ws.Range("I675").Value = 0
$ws.Range("I675").NumberFormat = "General"
$ws.Range("J675").Value = 221.9213591536121

# Row 676
$ws.Range("A676").NumberFormat = "@"
$ws.Range("A676").Value = "2024-09-03"
$ws.Range("A676").ClearFormats()
$ws.Range("C676").Value = 1865.599975585938
$ws.Range("D676").Value = 599.9400024414062
$ws.Range("E676").Value = 1114
$ws.Range("F676").Value = 178.4600067138672
$ws.Range("G676").Value = 659.0999755859375
$ws.Range("H676").Value = 22746.27993774414
$ws.Range("I676").Value = 0.003366991432787071
$ws.Range("I676").NumberFormat = "General"
$ws.Range("J676").Value = 222.6685664686348

# Row 677
$ws.Range("A677").NumberFormat = "@"
$ws.Range("A677").Value = "2024-09-04"
$ws.Range("A677").ClearFormats()
$ws.Range("C677").Value = 1871.900024414062
$ws.Range("D677").Value = 609
$ws.Range("E677").Value = 1127.900024414062
$ws.Range("F677").Value = 176.0200042724609
$ws.Range("G677").Value = 650.8499755859375
$ws.Range("H677").Value = 22782.94021606445
$ws.Range("I677").Value = 0.001611704349926693
$ws.Range("I677").NumberFormat = "General"
$ws.Range("J677").Value = 223.0274423658042

# Row 678
$ws.Range("A678").NumberFormat = "@"
$ws.Range("A678").Value = "2024-09-05"
$ws.Range("A678").ClearFormats()
$ws.Range("C678").Value = 1864.949951171875
$ws.Range("D678").Value = 602.1799926757812
$ws.Range("E678").Value = 1115.150024414062
$ws.Range("F678").Value = 173.4799957275391
$ws.Range("G678").Value = 643.8499755859375
$ws.Range("H678").Value = 22586.02963256836
$ws.Range("I678").Value = -0.008642896027846764
$ws.Range("I678").NumberFormat = "General"
$ws.Range("J678").Value = 221.09983937008

# Row 679
$ws.Range("A679").NumberFormat = "@"
$ws.Range("A679").Value = "2024-09-06"
$ws.Range("A679").ClearFormats()
$ws.Range("C679").Value = 1857.150024414062
$ws.Range("D679").Value = 597.2999877929688
$ws.Range("E679").Value = 1100
$ws.Range("F679").Value = 169.8500061035156
$ws.Range("G679").Value = 634.7000122070312
$ws.Range("H679").Value = 22350.45025634766
$ws.Range("I679").Value = -0.01043031378480992
$ws.Range("I679").NumberFormat = "General"
$ws.Range("J679").Value = 218.793698667679

# Row 680
$ws.Range("A680").NumberFormat = "@"
$ws.Range("A680").Value = "2024-09-09"
$ws.Range("A680").ClearFormats()
$ws.Range("C680").Value = 1860.449951171875
$ws.Range("D680").Value = 610.3400268554688
$ws.Range("E680").Value = 1104.150024414062
$ws.Range("F680").Value = 168.3300018310547
$ws.Range("G680").Value = 635.2000122070312
$ws.Range("H680").Value = 22404.27005004883
$ws.Range("I680").Value = 0.002407995950143633
$ws.Range("I680").NumberFormat = "General"
$ws.Range("J680").Value = 219.3205530079877

# Row 681
$ws.Range("A681").NumberFormat = "@"
$ws.Range("A681").Value = "2024-09-10"
$ws.Range("A681").ClearFormats()
$ws.Range("C681").Value = 1824.5
$ws.Range("D681").Value = 608
$ws.Range("E681").Value = 1113.199951171875
$ws.Range("F681").Value = 169.75
$ws.Range("G681").Value = 637.0499877929688
$ws.Range("H681").Value = 22289.99975585938
$ws.Range("I681").Value = -0.005100380147810443
$ws.Range("I681").NumberFormat = "General"
$ws.Range("J681").Value = 218.201934813419

# Row 682
$ws.Range("A682").NumberFormat = "@"
$ws.Range("A682").Value = "2024-09-11"
$ws.Range("A682").ClearFormats()
$ws.Range("C682").Value = 1833.150024414062
$ws.Range("D682").Value = 627.6599731445312
$ws.Range("E682").Value = 1112.599975585938
$ws.Range("F682").Value = 165.8800048828125
$ws.Range("G682").Value = 627.2000122070312
$ws.Range("H682").Value = 22284.95007324219
$ws.Range("I682").Value = -0.0002265447587481507
$ws.Range("I682").NumberFormat = "General"
$ws.Range("J682").Value = 218.1525023087383

# Row 683
$ws.Range("A683").NumberFormat = "@"
$ws.Range("A683").Value = "2024-09-12"
$ws.Range("A683").ClearFormats()
$ws.Range("C683").Value = 1854.849975585938
$ws.Range("D683").Value = 645.5999755859375
$ws.Range("E683").Value = 1120.099975585938
$ws.Range("F683").Value = 167.0200042724609
$ws.Range("G683").Value = 651.0999755859375
$ws.Range("H683").Value = 22615.88967895508
$ws.Range("I683").Value = 0.01485036334500268
$ws.Range("I683").NumberFormat = "General"
$ws.Range("J683").Value = 221.3921462326446

# Row 684
$ws.Range("A684").NumberFormat = "@"
$ws.Range("A684").Value = "2024-09-13"
$ws.Range("A684").ClearFormats()
$ws.Range("C684").Value = 1894.449951171875
$ws.Range("D684").Value = 646.6500244140625
$ws.Range("E684").Value = 1118.550048828125
$ws.Range("F684").Value = 167.25
$ws.Range("G684").Value = 633.4500122070312
$ws.Range("H684").Value = 22746.35009765625
$ws.Range("I684").Value = 0.005768529142701387
$ws.Range("I684").NumberFormat = "General"
$ws.Range("J684").Value = 222.6692532801529

# Row 685
$ws.Range("A685").NumberFormat = "@"
$ws.Range("A685").Value = "2024-09-16"
$ws.Range("A685").ClearFormats()
$ws.Range("C685").Value = 1857.599975585938
$ws.Range("D685").Value = 621.0499877929688
$ws.Range("E685").Value = 1115.849975585938
$ws.Range("F685").Value = 163.9600067138672
$ws.Range("G685").Value = 665.9500122070312
$ws.Range("H685").Value = 22506.51992797852
$ws.Range("I685").Value = -0.01054367705799297
$ws.Range("I685").NumberFormat = "General"
$ws.Range("J685").Value = 220.3215005828225

# Row 686
$ws.Range("A686").NumberFormat = "@"
$ws.Range("A686").Value = "2024-09-17"
$ws.Range("A686").ClearFormats()
$ws.Range("C686").Value = 1848.699951171875
$ws.Range("D686").Value = 649.6500244140625
$ws.Range("E686").Value = 1110.949951171875
$ws.Range("F686").Value = 160.6000061035156
$ws.Range("G686").Value = 666.3499755859375
$ws.Range("H686").Value = 22484.49969482422
$ws.Range("I686").Value = -0.0009783935155129372
$ws.Range("I686").NumberFormat = "General"
$ws.Range("J686").Value = 220.1059394553242

# Row 687
$ws.Range("A687").NumberFormat = "@"
$ws.Range("A687").Value = "2024-09-18"
$ws.Range("A687").ClearFormats()
$ws.Range("C687").Value = 1888.199951171875
$ws.Range("D687").Value = 646.7000122070312
$ws.Range("E687").Value = 1079.949951171875
$ws.Range("F687").Value = 158.5599975585938
$ws.Range("G687").Value = 651.7000122070312
$ws.Range("H687").Value = 22442.71960449219
$ws.Range("I687").Value = -0.00185817300358472
$ws.Range("I687").NumberFormat = "General"
$ws.Range("J687").Value = 219.6969445406996

# Row 688
$ws.Range("A688").NumberFormat = "@"
$ws.Range("A688").Value = "2024-09-19"
$ws.Range("A688").ClearFormats()
$ws.Range("C688").Value = 1890.400024414062
$ws.Range("D688").Value = 652.1500244140625
$ws.Range("E688").Value = 1054.449951171875
$ws.Range("F688").Value = 155.25
$ws.Range("G688").Value = 649.5999755859375
$ws.Range("H688").Value = 22292.29992675781
$ws.Range("I688").Value = -0.006702381903139165
$ws.Range("I688").NumberFormat = "General"
$ws.Range("J688").Value = 218.2244517154351

# Row 689
$ws.Range("A689").NumberFormat = "@"
$ws.Range("A689").Value = "2024-09-20"
$ws.Range("A689").ClearFormats()
$ws.Range("C689").Value = 1916.800048828125
$ws.Range("D689").Value = 654.4500122070312
$ws.Range("E689").Value = 1054.599975585938
$ws.Range("F689").Value = 161.4299926757812
$ws.Range("G689").Value = 665.1500244140625
$ws.Range("H689").Value = 22632.26013183594
$ws.Range("I689").Value = 0.01525011803156592
$ws.Range("I689").NumberFormat = "General"
$ws.Range("J689").Value = 221.5524003614692

# Row 690
$ws.Range("A690").NumberFormat = "@"
$ws.Range("A690").Value = "2024-09-23"
$ws.Range("A690").ClearFormats()
$ws.Range("C690").Value = 1919.949951171875
$ws.Range("D690").Value = 654.0999755859375
$ws.Range("E690").Value = 1055.25
$ws.Range("F690").Value = 159.5599975585938
$ws.Range("G690").Value = 672
$ws.Range("H690").Value = 22635.46960449219
$ws.Range("I690").Value = 0.0001418096397599883
$ws.Range("I690").NumberFormat = "General"
$ws.Range("J690").Value = 221.5838186275524

# Row 691
$ws.Range("A691").NumberFormat = "@"
$ws.Range("A691").Value = "2024-09-24"
$ws.Range("A691").ClearFormats()
$ws.Range("C691").Value = 1904.650024414062
$ws.Range("D691").Value = 646.8499755859375
$ws.Range("E691").Value = 1051.550048828125
$ws.Range("F691").Value = 158.7400054931641
$ws.Range("G691").Value = 675.25
$ws.Range("H691").Value = 22510.13034057617
$ws.Range("I691").Value = -0.005537294613544976
$ws.Range("I691").NumberFormat = "General"
$ws.Range("J691").Value = 220.3568437422173

# Row 692
$ws.Range("A692").NumberFormat = "@"
$ws.Range("A692").Value = "2024-09-25"
$ws.Range("A692").ClearFormats()
$ws.Range("C692").Value = 1928.5
$ws.Range("D692").Value = 633.2999877929688
$ws.Range("E692").Value = 1063.449951171875
$ws.Range("F692").Value = 156.9400024414062
$ws.Range("G692").Value = 667.3499755859375
$ws.Range("H692").Value = 22551.57971191406
$ws.Range("I692").Value = 0.001841365230265907
$ws.Range("I692").NumberFormat = "General"
$ws.Range("J692").Value = 220.7626011725354

# Row 693
$ws.Range("A693").NumberFormat = "@"
$ws.Range("A693").Value = "2024-09-26"
$ws.Range("A693").ClearFormats()
$ws.Range("C693").Value = 1982.800048828125
$ws.Range("D693").Value = 626.8499755859375
$ws.Range("E693").Value = 1068
$ws.Range("F693").Value = 156.8500061035156
$ws.Range("G693").Value = 665.3499755859375
$ws.Range("H693").Value = 22805.50018310547
$ws.Range("I693").Value = 0.01125954254358773
$ws.Range("I693").NumberFormat = "General"
$ws.Range("J693").Value = 223.2482870724707

# Row 694
$ws.Range("A694").NumberFormat = "@"
$ws.Range("A694").Value = "2024-09-27"
$ws.Range("A694").ClearFormats()
$ws.Range("C694").Value = 2010.699951171875
$ws.Range("D694").Value = 608.5499877929688
$ws.Range("E694").Value = 1075.949951171875
$ws.Range("F694").Value = 156.8099975585938
$ws.Range("G694").Value = 654.2999877929688
$ws.Range("H694").Value = 22858.51940917969
$ws.Range("I694").Value = 0.002324843816120109
$ws.Range("I694").NumberFormat = "General"
$ws.Range("J694").Value = 223.7673044721305
